# Sprint 41 test case report - fill in Day 3 results (API created & iOS executed)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Summary")

# Day 3 summary block (rows 15-17): Total testcase Written / Total Execution / Total Review
$ws.Range("C15").Value = 3614
$ws.Range("C16").Value = 1574
$ws.Range("C17").Value = 971

# Move the viewport/selection like the author left it: scrolled down with C17 selected
$window = $excel.ActiveWindow
$window.ScrollRow = 5
$window.ScrollColumn = 1
$window.Zoom = 100

$ws.Range("C17").Select() | Out-Null
